$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current (pre-edit) values for the columns that move: D,L,M,N,O,P,Q,R,S,T
$rows = @(2,3,5,6,7,8,9,10,11)
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

# Row data moves along this cycle: each row's values move to the NEXT row in
# the sequence (wrapping around at the end).
$chain = @(2,8,10,7,9,3,11,6,5)

for ($i = 0; $i -lt $chain.Length; $i++) {
    $src = $chain[$i]
    $dest = $chain[($i + 1) % $chain.Length]
    $data = $snapshot[$src]

    $ws.Cells.Item($dest, 4).Value2 = $data.D
    $ws.Cells.Item($dest, 12).Value2 = $data.L
    $ws.Cells.Item($dest, 13).Value2 = $data.M
    $ws.Cells.Item($dest, 14).Value2 = $data.N
    $ws.Cells.Item($dest, 15).Value2 = $data.O
    $ws.Cells.Item($dest, 16).Value2 = $data.P
    $ws.Cells.Item($dest, 17).Value2 = $data.Q
    $ws.Cells.Item($dest, 18).Value2 = $data.R
    $ws.Cells.Item($dest, 19).Value2 = $data.S
    $ws.Cells.Item($dest, 20).Value2 = $data.T
}
